$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "319.79"
Set-TextValue "E2" "-2.86%"
Set-TextValue "D3" "42.61"
Set-TextValue "E3" "-4.77%"
Set-TextValue "E4" "-5.55%"
Set-TextValue "D5" "0.08172"
Set-TextValue "E5" "-2.69%"
Set-TextValue "D6" "4.360"
Set-TextValue "D7" "1.763"
Set-TextValue "E7" "-12.62%"
Set-TextValue "D8" "0.9477"
Set-TextValue "E8" "-3.33%"
Set-TextValue "D9" "0.1124"
Set-TextValue "E9" "1.44%"
Set-TextValue "D10" "0.1876"
Set-TextValue "E10" "-1.80%"
Set-TextValue "D11" "0.04675"
Set-TextValue "E11" "-0.35%"
Set-TextValue "D12" "0.09363"
Set-TextValue "E12" "-3.01%"
Set-TextValue "D13" "7.496"
Set-TextValue "E13" "-20.97%"
Set-TextValue "D14" "0.1056"
Set-TextValue "E14" "-0.25%"
Set-TextValue "D15" "0.001289"
Set-TextValue "E15" "-1.72%"
Set-TextValue "D16" "0.005825"
Set-TextValue "E16" "-2.07%"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.351"
Set-TextValue "E17" "-1.07%"
Set-TextValue "B18" "BTSEToken"
Set-TextValue "C18" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D18" "2.532"
Set-TextValue "E18" "-0.19%"
Set-TextValue "B19" "BitpandaEcosystemToken"
Set-TextValue "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D19" "0.3362"
Set-TextValue "E19" "0.32%"
Set-TextValue "B20" "ProBitToken"
Set-TextValue "C20" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D20" "0.1347"
Set-TextValue "E20" "-0.87%"
Set-TextValue "B21" "ZBToken"
Set-TextValue "C21" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D21" "0.2557"
Set-TextValue "E21" "0.03%"
Set-TextValue "B22" "CoinExToken"
Set-TextValue "C22" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D22" "0.04180"
Set-TextValue "E22" "-0.06%"
Set-TextValue "B23" "BitKan"
Set-TextValue "C23" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D23" "0.001238"
Set-TextValue "E23" "-4.94%"
Set-TextValue "B24" "HotbitToken"
Set-TextValue "C24" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D24" "0.004280"
Set-TextValue "E24" "-3.50%"
Set-TextValue "E25" "-6.18%"
Set-TextValue "D26" "0.0002989"
Set-TextValue "E26" "0.06%"
Set-TextValue "E38" "-3.22%"
Set-TextValue "D39" "0.05615"
Set-TextValue "E39" "-0.65%"
Set-TextValue "D40" "0.008152"
Set-TextValue "E40" "4.28%"
Set-TextValue "D41" "0.1400"
Set-TextValue "E41" "-1.65%"
Set-TextValue "D42" "0.006550"
Set-TextValue "E42" "-11.31%"
Set-TextValue "D43" "0.002087"
Set-TextValue "E43" "-1.79%"
Set-TextValue "D44" "0.007618"
Set-TextValue "E44" "-11.78%"
Set-TextValue "D45" "0.3483"
Set-TextValue "E45" "2.96%"
Set-TextValue "D46" "0.00006796"
Set-TextValue "E46" "-1.29%"
Set-TextValue "D47" "0.00000000753"
Set-TextValue "E47" "-0.03%"
Set-TextValue "D48" "0.003355"
Set-TextValue "E48" "-3.67%"
Set-TextValue "D49" "0.004114"
Set-TextValue "E49" "16.12%"
Set-TextValue "E50" "-0.03%"
Set-TextValue "D51" "0.0002007"
Set-TextValue "E51" "-0.03%"
